$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    ("Play Dragon's Fire for Free - Review").
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.Text = "Meta description: Experience the fantasy world of Dragon's Fire slot. Stunning graphics, free spins, and a max win of 10,000x your bet. Play for free today."

# Bold just the "Meta description" label (first 16 characters) so it
# becomes its own run, matching the rest of the body text.
$labelRange = $d.Range($metaRange.Start, $metaRange.Start + 16)
$labelRange.Bold = 1

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Dragon's Fire for Free - Review"
#    paragraph near the end of the document. Search only the part of
#    the document AFTER the title/meta block we just built, so the
#    Find doesn't re-match the real title at the top.
# ---------------------------------------------------------------------
$oldTitleText = "Play Dragon's Fire for Free - Review"
$tailStart = $metaPara.Range.End

$titleFindRange = $d.Range($tailStart, $d.Content.End)
$titleFindRange.Find.MatchCase = $true
$dupFound = $titleFindRange.Find.Execute($oldTitleText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)

if ($dupFound) {
    # Extend by one character to also swallow the paragraph mark so the
    # whole paragraph (not just its text) is removed.
    $dupParaRange = $d.Range($titleFindRange.Start, $titleFindRange.End + 1)
    $dupParaRange.Delete()
}

# ---------------------------------------------------------------------
# 3) Replace the italic meta/prompt paragraph's text with the new image
#    prompt, keeping its italic run formatting intact.
# ---------------------------------------------------------------------
$oldPromptText = "Experience the fantasy world of Dragon's Fire slot. Stunning graphics, free spins, and a max win of 10,000x your bet. Play for free today."
$newPromptText = "Prompt: Create a feature image for Dragon's Fire that showcases a happy Maya warrior with glasses amidst the dragon-themed slot machine. The image should be in a cartoon style and must be eye-catching to suit the game's mesmerizing graphics."

$promptFindRange = $d.Range($tailStart, $d.Content.End)
$promptFindRange.Find.MatchCase = $true
$promptFound = $promptFindRange.Find.Execute($oldPromptText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)

if ($promptFound) {
    # Assign directly (instead of using Find's ReplaceWith) so Word's
    # smart-quote autocorrect doesn't mangle the straight apostrophes.
    $promptFindRange.Text = $newPromptText
}

Write-Host "Edits applied."
